$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"0.4767983333333334"
$ws.Cells.Item(2, 8).Value = [double]"1.430395"
$ws.Cells.Item(2, 9).Value = [double]"0.003723890400117776"
$ws.Cells.Item(2, 10).Value = [double]"0.003723890400117776"
$ws.Cells.Item(2, 11).Value = [double]"1"
$ws.Cells.Item(2, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 13).Value = [double]"0.007258333333333333"
$ws.Cells.Item(2, 14).Value = [double]"0.021775"
$ws.Cells.Item(2, 15).Value = [double]"0.000328667160253549"
$ws.Cells.Item(2, 16).Value = [double]"0.000328667160253549"
$ws.Cells.Item(2, 17).Value = [double]"0.003460761236111111"
$ws.Cells.Item(2, 18).Value = [double]"0.031146851125"
$ws.Cells.Item(2, 19).Value = [double]"1.223920482902162E-06"
$ws.Cells.Item(2, 20).Value = [double]"1.223920482902162E-06"

$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"0.4767983333333334"
$ws.Cells.Item(3, 8).Value = [double]"1.430395"
$ws.Cells.Item(3, 9).Value = [double]"0.003723890400117776"
$ws.Cells.Item(3, 10).Value = [double]"0.003723890400117776"
$ws.Cells.Item(3, 15).Value = [double]"0.7778551418094273"
$ws.Cells.Item(3, 16).Value = [double]"0.7778551418094272"
$ws.Cells.Item(3, 17).Value = [double]"8.190568598356668"
$ws.Cells.Item(3, 18).Value = [double]"73.71511738521001"
$ws.Cells.Item(3, 19).Value = [double]"0.002896647295266378"
$ws.Cells.Item(3, 20).Value = [double]"0.002896647295266377"

$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"0.4767983333333334"
$ws.Cells.Item(4, 8).Value = [double]"1.430395"
$ws.Cells.Item(4, 9).Value = [double]"0.003723890400117776"
$ws.Cells.Item(4, 10).Value = [double]"0.003723890400117776"
$ws.Cells.Item(4, 13).Value = [double]"4.898620999999999"
$ws.Cells.Item(4, 14).Value = [double]"14.695863"
$ws.Cells.Item(4, 15).Value = [double]"0.2218161910303192"
$ws.Cells.Item(4, 16).Value = [double]"0.2218161910303192"
$ws.Cells.Item(4, 17).Value = [double]"2.335654328431667"
$ws.Cells.Item(4, 18).Value = [double]"21.020888955885"
$ws.Cells.Item(4, 19).Value = [double]"0.0008260191843684966"
$ws.Cells.Item(4, 20).Value = [double]"0.0008260191843684965"

$ws.Cells.Item(5, 9).Value = [double]"0.8490200321922391"
$ws.Cells.Item(5, 10).Value = [double]"0.8490200321922391"
$ws.Cells.Item(5, 11).Value = [double]"1"
$ws.Cells.Item(5, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 13).Value = [double]"0.007258333333333333"
$ws.Cells.Item(5, 14).Value = [double]"0.021775"
$ws.Cells.Item(5, 15).Value = [double]"0.000328667160253549"
$ws.Cells.Item(5, 16).Value = [double]"0.000328667160253549"
$ws.Cells.Item(5, 17).Value = [double]"0.7890284891305555"
$ws.Cells.Item(5, 18).Value = [double]"7.101256402174999"
$ws.Cells.Item(5, 19).Value = [double]"0.000279045002979"
$ws.Cells.Item(5, 20).Value = [double]"0.000279045002979"

$ws.Cells.Item(6, 9).Value = [double]"0.8490200321922391"
$ws.Cells.Item(6, 10).Value = [double]"0.8490200321922391"
$ws.Cells.Item(6, 15).Value = [double]"0.7778551418094273"
$ws.Cells.Item(6, 16).Value = [double]"0.7778551418094272"
$ws.Cells.Item(6, 19).Value = [double]"0.6604145975399387"
$ws.Cells.Item(6, 20).Value = [double]"0.6604145975399386"

$ws.Cells.Item(7, 9).Value = [double]"0.8490200321922391"
$ws.Cells.Item(7, 10).Value = [double]"0.8490200321922391"
$ws.Cells.Item(7, 13).Value = [double]"4.898620999999999"
$ws.Cells.Item(7, 14).Value = [double]"14.695863"
$ws.Cells.Item(7, 15).Value = [double]"0.2218161910303192"
$ws.Cells.Item(7, 16).Value = [double]"0.2218161910303192"
$ws.Cells.Item(7, 17).Value = [double]"532.5122654126122"
$ws.Cells.Item(7, 18).Value = [double]"4792.61038871351"
$ws.Cells.Item(7, 19).Value = [double]"0.1883263896493215"
$ws.Cells.Item(7, 20).Value = [double]"0.1883263896493215"

$ws.Cells.Item(8, 7).Value = [double]"18.85432833333333"
$ws.Cells.Item(8, 8).Value = [double]"56.562985"
$ws.Cells.Item(8, 9).Value = [double]"0.1472560774076432"
$ws.Cells.Item(8, 10).Value = [double]"0.1472560774076432"
$ws.Cells.Item(8, 11).Value = [double]"1"
$ws.Cells.Item(8, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 13).Value = [double]"0.007258333333333333"
$ws.Cells.Item(8, 14).Value = [double]"0.021775"
$ws.Cells.Item(8, 15).Value = [double]"0.000328667160253549"
$ws.Cells.Item(8, 16).Value = [double]"0.000328667160253549"
$ws.Cells.Item(8, 17).Value = [double]"0.1368509998194444"
$ws.Cells.Item(8, 18).Value = [double]"1.231658998375"
$ws.Cells.Item(8, 19).Value = [double]"4.839823679164687E-05"
$ws.Cells.Item(8, 20).Value = [double]"4.839823679164687E-05"

$ws.Cells.Item(9, 7).Value = [double]"18.85432833333333"
$ws.Cells.Item(9, 8).Value = [double]"56.562985"
$ws.Cells.Item(9, 9).Value = [double]"0.1472560774076432"
$ws.Cells.Item(9, 10).Value = [double]"0.1472560774076432"
$ws.Cells.Item(9, 15).Value = [double]"0.7778551418094273"
$ws.Cells.Item(9, 16).Value = [double]"0.7778551418094272"
$ws.Cells.Item(9, 17).Value = [double]"323.8846673613367"
$ws.Cells.Item(9, 18).Value = [double]"2914.96200625203"
$ws.Cells.Item(9, 19).Value = [double]"0.1145438969742223"
$ws.Cells.Item(9, 20).Value = [double]"0.1145438969742223"

$ws.Cells.Item(10, 7).Value = [double]"18.85432833333333"
$ws.Cells.Item(10, 8).Value = [double]"56.562985"
$ws.Cells.Item(10, 9).Value = [double]"0.1472560774076432"
$ws.Cells.Item(10, 10).Value = [double]"0.1472560774076432"
$ws.Cells.Item(10, 13).Value = [double]"4.898620999999999"
$ws.Cells.Item(10, 14).Value = [double]"14.695863"
$ws.Cells.Item(10, 15).Value = [double]"0.2218161910303192"
$ws.Cells.Item(10, 16).Value = [double]"0.2218161910303192"
$ws.Cells.Item(10, 17).Value = [double]"92.36020871456165"
$ws.Cells.Item(10, 18).Value = [double]"831.2418784310549"
$ws.Cells.Item(10, 19).Value = [double]"0.03266378219662925"
$ws.Cells.Item(10, 20).Value = [double]"0.03266378219662926"
